$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "SI "
$ws2 = $wb.Worksheets.Item(2)   # "SI -erreur"

# --- Add the new "Feuil1" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$new.Name = "Feuil1"

# --- Move the "exercise notes" blocks that lived at the bottom of the two
#     data sheets onto the new sheet, each preceded by a little header ---

# Block coming from "SI -erreur" (rows 42:57) -> Feuil1 rows 2:17
$ws2.Range("A42:O57").Copy()
$new.Range("A2:O17").PasteSpecial()

# Block coming from "SI " (rows 42:52) -> Feuil1 rows 23:33
$ws1.Range("A42:O52").Copy()
$new.Range("A23:O33").PasteSpecial()

# Header labels — "SI sheet" must be created before "SI- erreur sheet" so
# the two new shared strings land in the same order as the target file.
$new.Range("A22").Value = "SI sheet"
$new.Range("A1").Value = "SI- erreur sheet"

# Remove the now-duplicated rows from the original sheets
$ws1.Rows("42:52").Delete()
$ws2.Rows("42:57").Delete()

# --- View bits ---
$ws2.Activate()
$excel.ActiveWindow.Zoom = 80
$ws2.Range("A42:Q60").Select() | Out-Null

$ws1.Activate()
$ws1.Range("A42:O58").Select() | Out-Null

$new.Activate()
$new.Range("G20").Select() | Out-Null
